# Daily attendance processing - 2025-12-18 01:28:45
# Rotate the "Recorded By" (column G) contributor list for each data row so
# that the last-listed contributor is moved to the front of the list -
# unless the last contributor already is "System", in which case the
# order is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $raw = $cell.Value2

    if ($raw -eq $null) { continue }

    $text = [string]$raw
    if ($text -eq "") { continue }

    $parts = $text -split ", "
    if ($parts.Count -gt 1 -and $parts[$parts.Count - 1] -cne "System") {
        $rotated = @($parts[$parts.Count - 1]) + $parts[0..($parts.Count - 2)]
        $newText = $rotated -join ", "
        $cell.Value = $newText
    }
}
